$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "gabriel"
$ws.Range("C2").Value = "araujo"
$ws.Range("D2").Value = "gabrielaraujo2334@gmail.com"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "31973413991"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = 1

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Eva"
$ws.Range("C3").Value = "Araujo"
$ws.Range("D3").Value = "EVAVILMA-ARAUJO@HOTMAIL.COM"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "31973413991"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = 2

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Manu"
$ws.Range("C4").Value = "Araujo"
$ws.Range("D4").Value = "adas@"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3141241"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = 3

# Apply the bold/border/centered header-style (style index used by row 1) to A2:A4
$ws.Range("B1").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
